{"js": "// Revert \"Creacion del metodo mul por escalar no miembro\":\n//  1) The bitacora row whose date cell reads \"27/01/2025\" goes back to \"26/01\".\n//  2) The whole following row (date \"28/01/2025\", the redimensionar/escalar entry) is removed.\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load the first-column (date) cell text for every row so we can find our targets by content.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.items[0].body.load(\"text\");\n}\nawait context.sync();\n\nlet dateRow = null;\nlet rowToDelete = null;\nfor (let i = 0; i < rows.items.length; i++) {\n  const text = rows.items[i].cells.items[0].body.text.trim();\n  if (text === \"27/01/2025\") {\n    dateRow = rows.items[i];\n    rowToDelete = rows.items[i + 1] || null;\n  }\n}\n\nif (!dateRow) {\n  throw new Error(\"Could not locate the bitacora row dated 27/01/2025\");\n}\n\n// Replace the (multi-run) \"27/01/2025\" text with a single \"26/01\" run.\nconst dateCell = dateRow.cells.items[0].body;\ndateCell.clear();\ndateCell.insertText(\"26/01\", Word.InsertLocation.start);\nawait context.sync();\n\n// Delete the next row entirely (the reverted \"28/01/2025\" entry).\nif (rowToDelete) {\n  rowToDelete.delete();\n  await context.sync();\n}\n", "ps1": "# Revert \"Creacion del metodo mul por escalar no miembro\":\n#  1) The bitacora row whose date cell reads \"27/01/2025\" goes back to \"26/01\".\n#  2) The whole following row (date \"28/01/2025\", the redimensionar/escalar entry) is removed.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Locate the row dated 27/01/2025 by scanning the first column.\n$targetIndex = 0\nfor ($i = 1; $i -le $table.Rows.Count; $i++) {\n    $cellText = $table.Rows.Item($i).Cells.Item(1).Range.Text\n    if ($cellText -match \"27/01/2025\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq 0) {\n    throw \"Could not locate the bitacora row dated 27/01/2025\"\n}\n\n# Replace the (multi-run) \"27/01/2025\" text with a single \"26/01\" run.\n# wdReplaceAll = 2 ; Wrap = wdFindStop(0) via Find scoped to this cell's range.\n$cellRange = $table.Rows.Item($targetIndex).Cells.Item(1).Range\n$cellRange.Find.Execute(\"27/01/2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"26/01\", 2)\n\n# Delete the next row entirely (the reverted \"28/01/2025\" entry).\n$table = $d.Tables.Item(1)\nif ($targetIndex -lt $table.Rows.Count) {\n    $table.Rows.Item($targetIndex + 1).Delete()\n}\n"}
